# Normalize the "Recorded By" column (column G) so that any comma-separated
# list of recorders that ends with "System" has "System" moved to the front
# of the list instead of the back, e.g.:
#   "dnasr281@gmail.com, System"               -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System"               -> "System, backup@backdoor.com"
#   "system, backup@backdoor.com, System"       -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = $ws.UsedRange.Row
$firstDataRow = $headerRow + 1
$lastRow = $headerRow + $ws.UsedRange.Rows.Count - 1

$colIndex = 7          # column G = "Recorded By"
$suffix = ", System"

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    $val = $cell.Value2

    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.EndsWith($suffix)) {
            $rest = $text.Substring(0, $text.Length - $suffix.Length)
            $cell.Value2 = "System, " + $rest
        }
    }
}
